$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "RECENTLYVISITED" access-control-flow block (mirrors the existing
# USER / TAG / RESOURCE / SUBSCRIPTION / SCHOOL blocks) in rows 26-29.
# ---------------------------------------------------------------------------

# 1) Copy the header-row formatting (row 21, e.g. SCHOOL's header) down to the
#    new header row 26.
$ws.Range("A21:J21").Copy()
$ws.Range("A26:J26").PasteSpecial(-4122)   # xlPasteFormats

# 2) Merge the route-name column for the new block FIRST (before formatting),
#    same as every other block - merging after formatting would re-split the
#    cell borders into top/middle/bottom variants that don't match the other
#    (pre-existing) merged blocks.
$ws.Range("B27:B29").Merge() | Out-Null

# 3) Copy the "all clear / no XOR-required cell" data-row formatting (row 24 -
#    SCHOOL's 3rd data row, fully green with no red required-field markers,
#    matching the new RECENTLYVISITED block which has no required pairing)
#    down across the three new data rows 27-29.
$ws.Range("A24:J24").Copy()
$ws.Range("A27:J29").PasteSpecial(-4122)   # xlPasteFormats

# 4) Populate the new header row with the same column headings used by every
#    other block.
$ws.Range("A26").Value = "ROLE"
$ws.Range("B26").Value = "ROUTE"
$ws.Range("C26").Value = "CREATE 1"
$ws.Range("D26").Value = "DELETE 1"
$ws.Range("E26").Value = "UPDATE 1"
$ws.Range("F26").Value = "READ 1"
$ws.Range("G26").Value = "BULK CREATE"
$ws.Range("H26").Value = "BULK DELETE"
$ws.Range("I26").Value = "BULK UPDATE"
$ws.Range("J26").Value = "BULK READ"

# 5) Populate the new data rows: STUDENT / TEACHER / ADMIN against the new
#    RECENTLYVISITED route.
$ws.Range("A27").Value = "STUDENT"
$ws.Range("B27").Value = "RECENTLYVISITED"
$ws.Range("A28").Value = "TEACHER"
$ws.Range("A29").Value = "ADMIN"

# ---------------------------------------------------------------------------
# Misc sheet-level cosmetic tweaks captured in the diff.
# ---------------------------------------------------------------------------

# Widen column B slightly to fit "RECENTLYVISITED".
$ws.Columns(2).ColumnWidth = 15.33

# Zoom in a bit and move the selection/active cell to A2.
$excel.ActiveWindow.Zoom = 108
$ws.Range("A2").Select() | Out-Null

Write-Host "RECENTLYVISITED block added"
